$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data set for rows 2-13 (title,firstname,lastname,age,dob,gender,officename,email,phone,address)
$data = @(
    @("ms", "kanathala", "saisree", 28, "1998-06-07", "female", "microsoft", "sree@gmail.com", 9865457667, "77-98/4,kpcolony,hyderbad,telangana,500032"),
    @("ms", "Soumyashree", "nand", 24, "1998-05-23", "female", "sreeclinic", "soumysree@gmail.com", 8786764322, "4-2/4,teacherscolony,nagarkunool,telangana,500023"),
    @("mr", "mothe", "anand", 45, "2002-07-23", "male", "accenture", "anand@gmail.com", 9897864348, "2-24-432,bhagirathacolony,jadcherla,telangana,500023"),
    @("ms", "budhuru", "shirisha", 24, "1998-09-26", "female", "innominds", "shirisha@gmail.com", 9550717673, "23,phpcolony,mahabubnagar,telangana,509893"),
    @("ms", "budhuru", "shirisha", 24, "1998-09-26", "female", "innominds", "shirisha@gmail.com", 9550717673, "23,phpcolony,mahabubnagar,telangana,509893"),
    @("ms", "Soumyashree", "nand", 24, "1998-05-23", "female", "sreeclinic", "soumysree@gmail.com", 8786764322, "4-2/4,teacherscolony,nagarkunool,telangana,500023"),
    @("mr", "mothe", "anand", 45, "2002-07-23", "male", "accenture", "anand@gmail.com", 9897864348, "2-24-432,bhagirathacolony,jadcherla,telangana,500023"),
    @("ms", "kanathala", "saisree", 28, "1998-06-07", "female", "microsoft", "sree@gmail.com", 9865457667, "77-98/4,kpcolony,hyderbad,telangana,500032"),
    @("ms", "kanathala", "saisree", 28, "1998-06-07", "female", "microsoft", "sree@gmail.com", 9865457667, "77-98/4,kpcolony,hyderbad,telangana,500032"),
    @("ms", "Soumyashree", "nand", 24, "1998-05-23", "female", "sreeclinic", "soumysree@gmail.com", 8786764322, "4-2/4,teacherscolony,nagarkunool,telangana,500023"),
    @("ms", "budhuru", "shirisha", 24, "1998-09-26", "female", "innominds", "shirisha@gmail.com", 9550717673, "23,phpcolony,mahabubnagar,telangana,509893"),
    @("mr", "mothe", "anand", 45, "2002-07-23", "male", "accenture", "anand@gmail.com", 9897864348, "2-24-432,bhagirathacolony,jadcherla,telangana,500023")
)

# Column E (dob) holds date-like strings (e.g. "1998-06-07") that must stay
# plain text rather than being auto-converted to date serials, so force a
# text number format on that column before writing into it.
$ws.Range("E2:E13").NumberFormat = "@"

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $ws.Cells.Item($row, 6).Value = $entry[5]
    $ws.Cells.Item($row, 7).Value = $entry[6]
    $ws.Cells.Item($row, 8).Value = $entry[7]
    $ws.Cells.Item($row, 9).Value = $entry[8]
    $ws.Cells.Item($row, 10).Value = $entry[9]
    $row = $row + 1
}
